$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: plate number changed
$ws.Range("B2").Value = "43C22665_C"

# A3 / B3: new label text + numeric value update
$ws.Range("A3").Value = "Tổng số phương tiện be"
$ws.Range("B3").Value = 71

# A4 / B4: new label text + value update (stored as text, not number)
$ws.Range("A4").Value = "Tổng số phương tiện fe"
$ws.Range("B4").Value = "'71"
$ws.Range("B4").Style = "Normal"

# Update the active selection to B4
$ws.Range("B4").Select()
